$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43
$ws.Cells.Item(43, 2).Value = 6286549
$ws.Cells.Item(43, 6).Value = 'Criciuma'
$ws.Cells.Item(43, 7).Value = 'CRB'
$ws.Cells.Item(43, 8).Value = 2
$ws.Cells.Item(43, 9).Value = 1
$ws.Cells.Item(43, 10).Value = 'H'
$ws.Cells.Item(43, 11).Value = 1.8
$ws.Cells.Item(43, 12).Value = 3.2
$ws.Cells.Item(43, 13).Value = 4.5
$ws.Cells.Item(43, 14).Value = 1.909
$ws.Cells.Item(43, 15).Value = 3.25
$ws.Cells.Item(43, 16).Value = 4.5
$ws.Cells.Item(43, 17).Value = -0.5
$ws.Cells.Item(43, 18).Value = 1.875
$ws.Cells.Item(43, 19).Value = 1.975
$ws.Cells.Item(43, 20).Value = 2
$ws.Cells.Item(43, 21).Value = 1.85
$ws.Cells.Item(43, 22).Value = 2
$ws.Cells.Item(43, 23).Value = 0.909
$ws.Cells.Item(43, 24).Value = -1
$ws.Cells.Item(43, 25).Value = -1
$ws.Cells.Item(43, 26).Value = 0.875
$ws.Cells.Item(43, 27).Value = -1
$ws.Cells.Item(43, 28).Value = 0.8500000000000001
$ws.Cells.Item(43, 29).Value = -1

# Row 44
$ws.Cells.Item(44, 2).Value = 6285530
$ws.Cells.Item(44, 6).Value = 'Chapecoense'
$ws.Cells.Item(44, 7).Value = 'Sport Recife'
$ws.Cells.Item(44, 8).Value = 1
$ws.Cells.Item(44, 9).Value = 1
$ws.Cells.Item(44, 10).Value = 'D'
$ws.Cells.Item(44, 11).Value = 3.6
$ws.Cells.Item(44, 12).Value = 3.1
$ws.Cells.Item(44, 13).Value = 2.05
$ws.Cells.Item(44, 14).Value = 3.5
$ws.Cells.Item(44, 15).Value = 3
$ws.Cells.Item(44, 16).Value = 2.3
$ws.Cells.Item(44, 17).Value = 0.25
$ws.Cells.Item(44, 18).Value = 1.875
$ws.Cells.Item(44, 19).Value = 1.925
$ws.Cells.Item(44, 20).Value = 1.75
$ws.Cells.Item(44, 21).Value = 1.775
$ws.Cells.Item(44, 22).Value = 2.025
$ws.Cells.Item(44, 23).Value = -1
$ws.Cells.Item(44, 24).Value = 2
$ws.Cells.Item(44, 25).Value = -1
$ws.Cells.Item(44, 26).Value = 0.4375
$ws.Cells.Item(44, 27).Value = -0.5
$ws.Cells.Item(44, 28).Value = 0.3875
$ws.Cells.Item(44, 29).Value = -0.5

# Row 49
$ws.Cells.Item(49, 2).Value = 6289129
$ws.Cells.Item(49, 6).Value = 'Guarani'
$ws.Cells.Item(49, 7).Value = 'Mirassol'
$ws.Cells.Item(49, 8).Value = 2
$ws.Cells.Item(49, 9).Value = 1
$ws.Cells.Item(49, 10).Value = 'H'
$ws.Cells.Item(49, 11).Value = 2.4
$ws.Cells.Item(49, 12).Value = 3.1
$ws.Cells.Item(49, 13).Value = 2.875
$ws.Cells.Item(49, 14).Value = 2.8
$ws.Cells.Item(49, 15).Value = 3.1
$ws.Cells.Item(49, 16).Value = 2.7
$ws.Cells.Item(49, 17).Value = 0
$ws.Cells.Item(49, 18).Value = 1.975
$ws.Cells.Item(49, 19).Value = 1.875
$ws.Cells.Item(49, 20).Value = 2
$ws.Cells.Item(49, 21).Value = 1.85
$ws.Cells.Item(49, 22).Value = 2
$ws.Cells.Item(49, 23).Value = 1.8
$ws.Cells.Item(49, 24).Value = -1
$ws.Cells.Item(49, 25).Value = -1
$ws.Cells.Item(49, 26).Value = 0.9750000000000001
$ws.Cells.Item(49, 27).Value = -1
$ws.Cells.Item(49, 28).Value = 0.8500000000000001
$ws.Cells.Item(49, 29).Value = -1

# Row 50
$ws.Cells.Item(50, 2).Value = 6281797
$ws.Cells.Item(50, 6).Value = 'Ituano'
$ws.Cells.Item(50, 7).Value = 'Ponte Preta'
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 9).Value = 1
$ws.Cells.Item(50, 10).Value = 'A'
$ws.Cells.Item(50, 11).Value = 2.3
$ws.Cells.Item(50, 12).Value = 3
$ws.Cells.Item(50, 13).Value = 3.2
$ws.Cells.Item(50, 14).Value = 2.05
$ws.Cells.Item(50, 15).Value = 3.1
$ws.Cells.Item(50, 16).Value = 4
$ws.Cells.Item(50, 17).Value = -0.25
$ws.Cells.Item(50, 18).Value = 1.75
$ws.Cells.Item(50, 19).Value = 2.05
$ws.Cells.Item(50, 20).Value = 2
$ws.Cells.Item(50, 21).Value = 1.9
$ws.Cells.Item(50, 22).Value = 1.9
$ws.Cells.Item(50, 23).Value = -1
$ws.Cells.Item(50, 24).Value = -1
$ws.Cells.Item(50, 25).Value = 3
$ws.Cells.Item(50, 26).Value = -1
$ws.Cells.Item(50, 27).Value = 1.05
$ws.Cells.Item(50, 28).Value = -1
$ws.Cells.Item(50, 29).Value = 0.8999999999999999

# Row 51
$ws.Cells.Item(51, 2).Value = 6282014
$ws.Cells.Item(51, 6).Value = 'Londrina'
$ws.Cells.Item(51, 7).Value = 'EC Juventude'
$ws.Cells.Item(51, 8).Value = 1
$ws.Cells.Item(51, 9).Value = 2
$ws.Cells.Item(51, 10).Value = 'A'
$ws.Cells.Item(51, 11).Value = 2.6
$ws.Cells.Item(51, 12).Value = 3.1
$ws.Cells.Item(51, 13).Value = 2.6
$ws.Cells.Item(51, 14).Value = 2.9
$ws.Cells.Item(51, 15).Value = 3
$ws.Cells.Item(51, 16).Value = 2.7
$ws.Cells.Item(51, 17).Value = 0
$ws.Cells.Item(51, 18).Value = 2
$ws.Cells.Item(51, 19).Value = 1.8
$ws.Cells.Item(51, 20).Value = 2
$ws.Cells.Item(51, 21).Value = 1.825
$ws.Cells.Item(51, 22).Value = 1.975
$ws.Cells.Item(51, 23).Value = -1
$ws.Cells.Item(51, 24).Value = -1
$ws.Cells.Item(51, 25).Value = 1.7
$ws.Cells.Item(51, 26).Value = -1
$ws.Cells.Item(51, 27).Value = 0.8
$ws.Cells.Item(51, 28).Value = 0.825
$ws.Cells.Item(51, 29).Value = -1

# Row 52
$ws.Cells.Item(52, 2).Value = 6281969
$ws.Cells.Item(52, 6).Value = 'Vitoria'
$ws.Cells.Item(52, 7).Value = 'Sampaio Correa'
$ws.Cells.Item(52, 8).Value = 2
$ws.Cells.Item(52, 9).Value = 1
$ws.Cells.Item(52, 10).Value = 'H'
$ws.Cells.Item(52, 11).Value = 1.615
$ws.Cells.Item(52, 12).Value = 3.6
$ws.Cells.Item(52, 13).Value = 5
$ws.Cells.Item(52, 14).Value = 1.75
$ws.Cells.Item(52, 15).Value = 3.3
$ws.Cells.Item(52, 16).Value = 5.5
$ws.Cells.Item(52, 17).Value = -0.75
$ws.Cells.Item(52, 18).Value = 1.975
$ws.Cells.Item(52, 19).Value = 1.825
$ws.Cells.Item(52, 20).Value = 2
$ws.Cells.Item(52, 21).Value = 1.825
$ws.Cells.Item(52, 22).Value = 1.975
$ws.Cells.Item(52, 23).Value = 0.75
$ws.Cells.Item(52, 24).Value = -1
$ws.Cells.Item(52, 25).Value = -1
$ws.Cells.Item(52, 26).Value = 0.4875
$ws.Cells.Item(52, 27).Value = -0.5
$ws.Cells.Item(52, 28).Value = 0.825
$ws.Cells.Item(52, 29).Value = -1

# Row 88
$ws.Cells.Item(88, 2).Value = 6281807
$ws.Cells.Item(88, 6).Value = 'ABC'
$ws.Cells.Item(88, 7).Value = 'Guarani'
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 1
$ws.Cells.Item(88, 10).Value = 'A'
$ws.Cells.Item(88, 11).Value = 2.75
$ws.Cells.Item(88, 12).Value = 2.875
$ws.Cells.Item(88, 13).Value = 2.55
$ws.Cells.Item(88, 14).Value = 2.5
$ws.Cells.Item(88, 15).Value = 2.9
$ws.Cells.Item(88, 16).Value = 3.25
$ws.Cells.Item(88, 17).Value = -0.25
$ws.Cells.Item(88, 18).Value = 2.075
$ws.Cells.Item(88, 19).Value = 1.725
$ws.Cells.Item(88, 20).Value = 2
$ws.Cells.Item(88, 21).Value = 2
$ws.Cells.Item(88, 22).Value = 1.8
$ws.Cells.Item(88, 23).Value = -1
$ws.Cells.Item(88, 24).Value = -1
$ws.Cells.Item(88, 25).Value = 2.25
$ws.Cells.Item(88, 26).Value = -1
$ws.Cells.Item(88, 27).Value = 0.7250000000000001
$ws.Cells.Item(88, 28).Value = -1
$ws.Cells.Item(88, 29).Value = 0.8

# Row 89
$ws.Cells.Item(89, 2).Value = 6285538
$ws.Cells.Item(89, 6).Value = 'Sport Recife'
$ws.Cells.Item(89, 7).Value = 'Vitoria'
$ws.Cells.Item(89, 8).Value = 1
$ws.Cells.Item(89, 9).Value = 2
$ws.Cells.Item(89, 10).Value = 'A'
$ws.Cells.Item(89, 11).Value = 1.666
$ws.Cells.Item(89, 12).Value = 3.3
$ws.Cells.Item(89, 13).Value = 5
$ws.Cells.Item(89, 14).Value = 1.7
$ws.Cells.Item(89, 15).Value = 3.4
$ws.Cells.Item(89, 16).Value = 6
$ws.Cells.Item(89, 17).Value = -0.75
$ws.Cells.Item(89, 18).Value = 1.925
$ws.Cells.Item(89, 19).Value = 1.875
$ws.Cells.Item(89, 20).Value = 2
$ws.Cells.Item(89, 21).Value = 1.825
$ws.Cells.Item(89, 22).Value = 1.975
$ws.Cells.Item(89, 23).Value = -1
$ws.Cells.Item(89, 24).Value = -1
$ws.Cells.Item(89, 25).Value = 5
$ws.Cells.Item(89, 26).Value = -1
$ws.Cells.Item(89, 27).Value = 0.875
$ws.Cells.Item(89, 28).Value = 0.825
$ws.Cells.Item(89, 29).Value = -1

# Row 90
$ws.Cells.Item(90, 2).Value = 6282102
$ws.Cells.Item(90, 6).Value = 'Ceara'
$ws.Cells.Item(90, 7).Value = 'Vila Nova'
$ws.Cells.Item(90, 8).Value = 1
$ws.Cells.Item(90, 9).Value = 0
$ws.Cells.Item(90, 10).Value = 'H'
$ws.Cells.Item(90, 11).Value = 2.4
$ws.Cells.Item(90, 12).Value = 2.875
$ws.Cells.Item(90, 13).Value = 3
$ws.Cells.Item(90, 14).Value = 2.7
$ws.Cells.Item(90, 15).Value = 2.8
$ws.Cells.Item(90, 16).Value = 3.1
$ws.Cells.Item(90, 17).Value = 0
$ws.Cells.Item(90, 18).Value = 1.7
$ws.Cells.Item(90, 19).Value = 2.2
$ws.Cells.Item(90, 20).Value = 1.75
$ws.Cells.Item(90, 21).Value = 1.85
$ws.Cells.Item(90, 22).Value = 2
$ws.Cells.Item(90, 23).Value = 1.7
$ws.Cells.Item(90, 24).Value = -1
$ws.Cells.Item(90, 25).Value = -1
$ws.Cells.Item(90, 26).Value = 0.7
$ws.Cells.Item(90, 27).Value = -1
$ws.Cells.Item(90, 28).Value = -1
$ws.Cells.Item(90, 29).Value = 1

# Row 94
$ws.Cells.Item(94, 2).Value = 6285669
$ws.Cells.Item(94, 6).Value = 'Vila Nova'
$ws.Cells.Item(94, 7).Value = 'ABC'
$ws.Cells.Item(94, 8).Value = 1
$ws.Cells.Item(94, 9).Value = 1
$ws.Cells.Item(94, 10).Value = 'D'
$ws.Cells.Item(94, 11).Value = 1.4
$ws.Cells.Item(94, 12).Value = 4
$ws.Cells.Item(94, 13).Value = 6.5
$ws.Cells.Item(94, 14).Value = 1.5
$ws.Cells.Item(94, 15).Value = 3.8
$ws.Cells.Item(94, 16).Value = 8.5
$ws.Cells.Item(94, 17).Value = -1
$ws.Cells.Item(94, 18).Value = 1.85
$ws.Cells.Item(94, 19).Value = 1.95
$ws.Cells.Item(94, 20).Value = 2
$ws.Cells.Item(94, 21).Value = 1.85
$ws.Cells.Item(94, 22).Value = 1.95
$ws.Cells.Item(94, 23).Value = -1
$ws.Cells.Item(94, 24).Value = 2.8
$ws.Cells.Item(94, 25).Value = -1
$ws.Cells.Item(94, 26).Value = -1
$ws.Cells.Item(94, 27).Value = 0.95
$ws.Cells.Item(94, 28).Value = 0
$ws.Cells.Item(94, 29).Value = -0

# Row 95
$ws.Cells.Item(95, 2).Value = 6285539
$ws.Cells.Item(95, 6).Value = 'Criciuma'
$ws.Cells.Item(95, 7).Value = 'Gremio Novorizontino'
$ws.Cells.Item(95, 8).Value = 0
$ws.Cells.Item(95, 9).Value = 1
$ws.Cells.Item(95, 10).Value = 'A'
$ws.Cells.Item(95, 11).Value = 2.3
$ws.Cells.Item(95, 12).Value = 2.875
$ws.Cells.Item(95, 13).Value = 3
$ws.Cells.Item(95, 14).Value = 2.45
$ws.Cells.Item(95, 15).Value = 2.9
$ws.Cells.Item(95, 16).Value = 3.4
$ws.Cells.Item(95, 17).Value = -0.25
$ws.Cells.Item(95, 18).Value = 2.025
$ws.Cells.Item(95, 19).Value = 1.775
$ws.Cells.Item(95, 20).Value = 2
$ws.Cells.Item(95, 21).Value = 2
$ws.Cells.Item(95, 22).Value = 1.8
$ws.Cells.Item(95, 23).Value = -1
$ws.Cells.Item(95, 24).Value = -1
$ws.Cells.Item(95, 25).Value = 2.4
$ws.Cells.Item(95, 26).Value = -1
$ws.Cells.Item(95, 27).Value = 0.7749999999999999
$ws.Cells.Item(95, 28).Value = -1
$ws.Cells.Item(95, 29).Value = 0.8

# Row 100
$ws.Cells.Item(100, 2).Value = 6286300
$ws.Cells.Item(100, 6).Value = 'Vitoria'
$ws.Cells.Item(100, 7).Value = 'Chapecoense'
$ws.Cells.Item(100, 8).Value = 1
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 10).Value = 'H'
$ws.Cells.Item(100, 11).Value = 1.7
$ws.Cells.Item(100, 12).Value = 3.25
$ws.Cells.Item(100, 13).Value = 4.5
$ws.Cells.Item(100, 14).Value = 1.7
$ws.Cells.Item(100, 15).Value = 3.5
$ws.Cells.Item(100, 16).Value = 5.75
$ws.Cells.Item(100, 17).Value = -0.75
$ws.Cells.Item(100, 18).Value = 1.9
$ws.Cells.Item(100, 19).Value = 1.95
$ws.Cells.Item(100, 20).Value = 2.25
$ws.Cells.Item(100, 21).Value = 2.05
$ws.Cells.Item(100, 22).Value = 1.8
$ws.Cells.Item(100, 23).Value = 0.7
$ws.Cells.Item(100, 24).Value = -1
$ws.Cells.Item(100, 25).Value = -1
$ws.Cells.Item(100, 26).Value = 0.45
$ws.Cells.Item(100, 27).Value = -0.5
$ws.Cells.Item(100, 28).Value = -1
$ws.Cells.Item(100, 29).Value = 0.8

# Row 101
$ws.Cells.Item(101, 2).Value = 6281811
$ws.Cells.Item(101, 6).Value = 'Londrina'
$ws.Cells.Item(101, 7).Value = 'Botafogo SP'
$ws.Cells.Item(101, 8).Value = 1
$ws.Cells.Item(101, 9).Value = 2
$ws.Cells.Item(101, 10).Value = 'A'
$ws.Cells.Item(101, 11).Value = 2.5
$ws.Cells.Item(101, 12).Value = 3
$ws.Cells.Item(101, 13).Value = 2.625
$ws.Cells.Item(101, 14).Value = 2.6
$ws.Cells.Item(101, 15).Value = 3
$ws.Cells.Item(101, 16).Value = 2.9
$ws.Cells.Item(101, 17).Value = 0
$ws.Cells.Item(101, 18).Value = 1.775
$ws.Cells.Item(101, 19).Value = 2.025
$ws.Cells.Item(101, 20).Value = 1.75
$ws.Cells.Item(101, 21).Value = 1.8
$ws.Cells.Item(101, 22).Value = 2
$ws.Cells.Item(101, 23).Value = -1
$ws.Cells.Item(101, 24).Value = -1
$ws.Cells.Item(101, 25).Value = 1.9
$ws.Cells.Item(101, 26).Value = -1
$ws.Cells.Item(101, 27).Value = 1.025
$ws.Cells.Item(101, 28).Value = 0.8
$ws.Cells.Item(101, 29).Value = -1

# Row 115
$ws.Cells.Item(115, 2).Value = 6285545
$ws.Cells.Item(115, 6).Value = 'Criciuma'
$ws.Cells.Item(115, 7).Value = 'Ponte Preta'
$ws.Cells.Item(115, 8).Value = 2
$ws.Cells.Item(115, 9).Value = 1
$ws.Cells.Item(115, 10).Value = 'H'
$ws.Cells.Item(115, 11).Value = 1.727
$ws.Cells.Item(115, 12).Value = 3.2
$ws.Cells.Item(115, 13).Value = 4.5
$ws.Cells.Item(115, 14).Value = 1.75
$ws.Cells.Item(115, 15).Value = 3.5
$ws.Cells.Item(115, 16).Value = 5.25
$ws.Cells.Item(115, 17).Value = -0.75
$ws.Cells.Item(115, 18).Value = 2
$ws.Cells.Item(115, 19).Value = 1.8
$ws.Cells.Item(115, 20).Value = 2
$ws.Cells.Item(115, 21).Value = 1.875
$ws.Cells.Item(115, 22).Value = 1.925
$ws.Cells.Item(115, 23).Value = 0.75
$ws.Cells.Item(115, 24).Value = -1
$ws.Cells.Item(115, 25).Value = -1
$ws.Cells.Item(115, 26).Value = 0.5
$ws.Cells.Item(115, 27).Value = -0.5
$ws.Cells.Item(115, 28).Value = 0.875
$ws.Cells.Item(115, 29).Value = -1

# Row 116
$ws.Cells.Item(116, 2).Value = 6281816
$ws.Cells.Item(116, 6).Value = 'Ituano'
$ws.Cells.Item(116, 7).Value = 'Tombense MG'
$ws.Cells.Item(116, 8).Value = 1
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 10).Value = 'H'
$ws.Cells.Item(116, 11).Value = 2.1
$ws.Cells.Item(116, 12).Value = 3
$ws.Cells.Item(116, 13).Value = 3.25
$ws.Cells.Item(116, 14).Value = 1.95
$ws.Cells.Item(116, 15).Value = 3.2
$ws.Cells.Item(116, 16).Value = 4.5
$ws.Cells.Item(116, 17).Value = -0.5
$ws.Cells.Item(116, 18).Value = 1.975
$ws.Cells.Item(116, 19).Value = 1.825
$ws.Cells.Item(116, 20).Value = 2
$ws.Cells.Item(116, 21).Value = 1.95
$ws.Cells.Item(116, 22).Value = 1.85
$ws.Cells.Item(116, 23).Value = 0.95
$ws.Cells.Item(116, 24).Value = -1
$ws.Cells.Item(116, 25).Value = -1
$ws.Cells.Item(116, 26).Value = 0.9750000000000001
$ws.Cells.Item(116, 27).Value = -1
$ws.Cells.Item(116, 28).Value = -1
$ws.Cells.Item(116, 29).Value = 0.8500000000000001

# Row 117
$ws.Cells.Item(117, 2).Value = 6285671
$ws.Cells.Item(117, 6).Value = 'Vitoria'
$ws.Cells.Item(117, 7).Value = 'ABC'
$ws.Cells.Item(117, 8).Value = 2
$ws.Cells.Item(117, 9).Value = 0
$ws.Cells.Item(117, 10).Value = 'H'
$ws.Cells.Item(117, 11).Value = 1.4
$ws.Cells.Item(117, 12).Value = 4
$ws.Cells.Item(117, 13).Value = 6.5
$ws.Cells.Item(117, 14).Value = 1.5
$ws.Cells.Item(117, 15).Value = 4
$ws.Cells.Item(117, 16).Value = 7.5
$ws.Cells.Item(117, 17).Value = -1
$ws.Cells.Item(117, 18).Value = 1.875
$ws.Cells.Item(117, 19).Value = 1.925
$ws.Cells.Item(117, 20).Value = 2
$ws.Cells.Item(117, 21).Value = 1.85
$ws.Cells.Item(117, 22).Value = 1.95
$ws.Cells.Item(117, 23).Value = 0.5
$ws.Cells.Item(117, 24).Value = -1
$ws.Cells.Item(117, 25).Value = -1
$ws.Cells.Item(117, 26).Value = 0.875
$ws.Cells.Item(117, 27).Value = -1
$ws.Cells.Item(117, 28).Value = 0
$ws.Cells.Item(117, 29).Value = -0

# Row 118
$ws.Cells.Item(118, 2).Value = 6282019
$ws.Cells.Item(118, 6).Value = 'EC Juventude'
$ws.Cells.Item(118, 7).Value = 'Gremio Novorizontino'
$ws.Cells.Item(118, 8).Value = 1
$ws.Cells.Item(118, 9).Value = 0
$ws.Cells.Item(118, 10).Value = 'H'
$ws.Cells.Item(118, 11).Value = 2.375
$ws.Cells.Item(118, 12).Value = 3
$ws.Cells.Item(118, 13).Value = 2.75
$ws.Cells.Item(118, 14).Value = 2.375
$ws.Cells.Item(118, 15).Value = 3.1
$ws.Cells.Item(118, 16).Value = 3.3
$ws.Cells.Item(118, 17).Value = -0.25
$ws.Cells.Item(118, 18).Value = 2.05
$ws.Cells.Item(118, 19).Value = 1.8
$ws.Cells.Item(118, 20).Value = 1.75
$ws.Cells.Item(118, 21).Value = 1.775
$ws.Cells.Item(118, 22).Value = 2.1
$ws.Cells.Item(118, 23).Value = 1.375
$ws.Cells.Item(118, 24).Value = -1
$ws.Cells.Item(118, 25).Value = -1
$ws.Cells.Item(118, 26).Value = 1.05
$ws.Cells.Item(118, 27).Value = -1
$ws.Cells.Item(118, 28).Value = -1
$ws.Cells.Item(118, 29).Value = 1.1

# Row 119
$ws.Cells.Item(119, 2).Value = 6287040
$ws.Cells.Item(119, 6).Value = 'Londrina'
$ws.Cells.Item(119, 7).Value = 'Chapecoense'
$ws.Cells.Item(119, 8).Value = 1
$ws.Cells.Item(119, 9).Value = 1
$ws.Cells.Item(119, 10).Value = 'D'
$ws.Cells.Item(119, 11).Value = 2.4
$ws.Cells.Item(119, 12).Value = 2.875
$ws.Cells.Item(119, 13).Value = 2.875
$ws.Cells.Item(119, 14).Value = 2.5
$ws.Cells.Item(119, 15).Value = 2.875
$ws.Cells.Item(119, 16).Value = 3.2
$ws.Cells.Item(119, 17).Value = -0.25
$ws.Cells.Item(119, 18).Value = 2.1
$ws.Cells.Item(119, 19).Value = 1.775
$ws.Cells.Item(119, 20).Value = 2
$ws.Cells.Item(119, 21).Value = 2.05
$ws.Cells.Item(119, 22).Value = 1.8
$ws.Cells.Item(119, 23).Value = -1
$ws.Cells.Item(119, 24).Value = 1.875
$ws.Cells.Item(119, 25).Value = -1
$ws.Cells.Item(119, 26).Value = -0.5
$ws.Cells.Item(119, 27).Value = 0.3875
$ws.Cells.Item(119, 28).Value = 0
$ws.Cells.Item(119, 29).Value = -0

# Row 120
$ws.Cells.Item(120, 2).Value = 6281815
$ws.Cells.Item(120, 6).Value = 'Guarani'
$ws.Cells.Item(120, 7).Value = 'Ceara'
$ws.Cells.Item(120, 8).Value = 0
$ws.Cells.Item(120, 9).Value = 0
$ws.Cells.Item(120, 10).Value = 'D'
$ws.Cells.Item(120, 11).Value = 1.95
$ws.Cells.Item(120, 12).Value = 3.25
$ws.Cells.Item(120, 13).Value = 3.4
$ws.Cells.Item(120, 14).Value = 2.05
$ws.Cells.Item(120, 15).Value = 3.25
$ws.Cells.Item(120, 16).Value = 4
$ws.Cells.Item(120, 17).Value = -0.5
$ws.Cells.Item(120, 18).Value = 2.025
$ws.Cells.Item(120, 19).Value = 1.775
$ws.Cells.Item(120, 20).Value = 2
$ws.Cells.Item(120, 21).Value = 1.95
$ws.Cells.Item(120, 22).Value = 1.85
$ws.Cells.Item(120, 23).Value = -1
$ws.Cells.Item(120, 24).Value = 2.25
$ws.Cells.Item(120, 25).Value = -1
$ws.Cells.Item(120, 26).Value = -1
$ws.Cells.Item(120, 27).Value = 0.7749999999999999
$ws.Cells.Item(120, 28).Value = -1
$ws.Cells.Item(120, 29).Value = 0.8500000000000001

# Row 121
$ws.Cells.Item(121, 2).Value = 6285544
$ws.Cells.Item(121, 6).Value = 'Vila Nova'
$ws.Cells.Item(121, 7).Value = 'Sport Recife'
$ws.Cells.Item(121, 8).Value = 0
$ws.Cells.Item(121, 9).Value = 1
$ws.Cells.Item(121, 10).Value = 'A'
$ws.Cells.Item(121, 11).Value = 2.2
$ws.Cells.Item(121, 12).Value = 2.875
$ws.Cells.Item(121, 13).Value = 3.2
$ws.Cells.Item(121, 14).Value = 2.25
$ws.Cells.Item(121, 15).Value = 2.875
$ws.Cells.Item(121, 16).Value = 3.8
$ws.Cells.Item(121, 17).Value = -0.25
$ws.Cells.Item(121, 18).Value = 1.9
$ws.Cells.Item(121, 19).Value = 1.95
$ws.Cells.Item(121, 20).Value = 1.75
$ws.Cells.Item(121, 21).Value = 1.825
$ws.Cells.Item(121, 22).Value = 2.025
$ws.Cells.Item(121, 23).Value = -1
$ws.Cells.Item(121, 24).Value = -1
$ws.Cells.Item(121, 25).Value = 2.8
$ws.Cells.Item(121, 26).Value = -1
$ws.Cells.Item(121, 27).Value = 0.95
$ws.Cells.Item(121, 28).Value = -1
$ws.Cells.Item(121, 29).Value = 1.025

# Row 130
$ws.Cells.Item(130, 2).Value = 6959080
$ws.Cells.Item(130, 6).Value = 'Atletico GO'
$ws.Cells.Item(130, 7).Value = 'Tombense MG'
$ws.Cells.Item(130, 8).Value = 3
$ws.Cells.Item(130, 9).Value = 2
$ws.Cells.Item(130, 10).Value = 'H'
$ws.Cells.Item(130, 11).Value = 1.7
$ws.Cells.Item(130, 12).Value = 3.25
$ws.Cells.Item(130, 13).Value = 4.5
$ws.Cells.Item(130, 14).Value = 1.727
$ws.Cells.Item(130, 15).Value = 3.6
$ws.Cells.Item(130, 16).Value = 5.25
$ws.Cells.Item(130, 17).Value = -0.75
$ws.Cells.Item(130, 18).Value = 1.975
$ws.Cells.Item(130, 19).Value = 1.825
$ws.Cells.Item(130, 20).Value = 2.25
$ws.Cells.Item(130, 21).Value = 1.975
$ws.Cells.Item(130, 22).Value = 1.825
$ws.Cells.Item(130, 23).Value = 0.7270000000000001
$ws.Cells.Item(130, 24).Value = -1
$ws.Cells.Item(130, 25).Value = -1
$ws.Cells.Item(130, 26).Value = 0.4875
$ws.Cells.Item(130, 27).Value = -0.5
$ws.Cells.Item(130, 28).Value = 0.9750000000000001
$ws.Cells.Item(130, 29).Value = -1

# Row 131
$ws.Cells.Item(131, 2).Value = 6286301
$ws.Cells.Item(131, 6).Value = 'Londrina'
$ws.Cells.Item(131, 7).Value = 'Vitoria'
$ws.Cells.Item(131, 8).Value = 2
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 'H'
$ws.Cells.Item(131, 11).Value = 3.1
$ws.Cells.Item(131, 12).Value = 3.1
$ws.Cells.Item(131, 13).Value = 2.15
$ws.Cells.Item(131, 14).Value = 3
$ws.Cells.Item(131, 15).Value = 3.2
$ws.Cells.Item(131, 16).Value = 2.45
$ws.Cells.Item(131, 17).Value = 0.25
$ws.Cells.Item(131, 18).Value = 1.725
$ws.Cells.Item(131, 19).Value = 2.075
$ws.Cells.Item(131, 20).Value = 1.75
$ws.Cells.Item(131, 21).Value = 1.775
$ws.Cells.Item(131, 22).Value = 2.025
$ws.Cells.Item(131, 23).Value = 2
$ws.Cells.Item(131, 24).Value = -1
$ws.Cells.Item(131, 25).Value = -1
$ws.Cells.Item(131, 26).Value = 0.7250000000000001
$ws.Cells.Item(131, 27).Value = -1
$ws.Cells.Item(131, 28).Value = 0.3875
$ws.Cells.Item(131, 29).Value = -0.5

# Row 188
$ws.Cells.Item(188, 2).Value = 6330569
$ws.Cells.Item(188, 6).Value = 'Sampaio Correa'
$ws.Cells.Item(188, 7).Value = 'Chapecoense'
$ws.Cells.Item(188, 8).Value = 2
$ws.Cells.Item(188, 9).Value = 0
$ws.Cells.Item(188, 10).Value = 'H'
$ws.Cells.Item(188, 11).Value = 2.05
$ws.Cells.Item(188, 12).Value = 2.9
$ws.Cells.Item(188, 13).Value = 4
$ws.Cells.Item(188, 14).Value = 2.05
$ws.Cells.Item(188, 15).Value = 2.9
$ws.Cells.Item(188, 16).Value = 4.5
$ws.Cells.Item(188, 17).Value = -0.5
$ws.Cells.Item(188, 18).Value = 2.05
$ws.Cells.Item(188, 19).Value = 1.75
$ws.Cells.Item(188, 20).Value = 1.75
$ws.Cells.Item(188, 21).Value = 1.775
$ws.Cells.Item(188, 22).Value = 2.025
$ws.Cells.Item(188, 23).Value = 1.05
$ws.Cells.Item(188, 24).Value = -1
$ws.Cells.Item(188, 25).Value = -1
$ws.Cells.Item(188, 26).Value = 1.05
$ws.Cells.Item(188, 27).Value = -1
$ws.Cells.Item(188, 28).Value = 0.3875
$ws.Cells.Item(188, 29).Value = -0.5

# Row 189
$ws.Cells.Item(189, 2).Value = 6285765
$ws.Cells.Item(189, 6).Value = 'Botafogo SP'
$ws.Cells.Item(189, 7).Value = 'Atletico GO'
$ws.Cells.Item(189, 8).Value = 1
$ws.Cells.Item(189, 9).Value = 0
$ws.Cells.Item(189, 10).Value = 'H'
$ws.Cells.Item(189, 11).Value = 2.875
$ws.Cells.Item(189, 12).Value = 3
$ws.Cells.Item(189, 13).Value = 2.55
$ws.Cells.Item(189, 14).Value = 3.3
$ws.Cells.Item(189, 15).Value = 2.9
$ws.Cells.Item(189, 16).Value = 2.4
$ws.Cells.Item(189, 17).Value = 0.25
$ws.Cells.Item(189, 18).Value = 1.775
$ws.Cells.Item(189, 19).Value = 2.025
$ws.Cells.Item(189, 20).Value = 1.75
$ws.Cells.Item(189, 21).Value = 1.825
$ws.Cells.Item(189, 22).Value = 1.975
$ws.Cells.Item(189, 23).Value = 2.3
$ws.Cells.Item(189, 24).Value = -1
$ws.Cells.Item(189, 25).Value = -1
$ws.Cells.Item(189, 26).Value = 0.7749999999999999
$ws.Cells.Item(189, 27).Value = -1
$ws.Cells.Item(189, 28).Value = -1
$ws.Cells.Item(189, 29).Value = 0.9750000000000001

# Row 216
$ws.Cells.Item(216, 2).Value = 6537870
$ws.Cells.Item(216, 6).Value = 'Atletico GO'
$ws.Cells.Item(216, 7).Value = 'Ituano'
$ws.Cells.Item(216, 8).Value = 1
$ws.Cells.Item(216, 9).Value = 0
$ws.Cells.Item(216, 10).Value = 'H'
$ws.Cells.Item(216, 11).Value = 1.615
$ws.Cells.Item(216, 12).Value = 3.5
$ws.Cells.Item(216, 13).Value = 5.25
$ws.Cells.Item(216, 14).Value = 1.65
$ws.Cells.Item(216, 15).Value = 3.6
$ws.Cells.Item(216, 16).Value = 5.75
$ws.Cells.Item(216, 17).Value = -0.75
$ws.Cells.Item(216, 18).Value = 1.8
$ws.Cells.Item(216, 19).Value = 2.05
$ws.Cells.Item(216, 20).Value = 2.25
$ws.Cells.Item(216, 21).Value = 1.9
$ws.Cells.Item(216, 22).Value = 1.95
$ws.Cells.Item(216, 23).Value = 0.6499999999999999
$ws.Cells.Item(216, 24).Value = -1
$ws.Cells.Item(216, 25).Value = -1
$ws.Cells.Item(216, 26).Value = 0.4
$ws.Cells.Item(216, 27).Value = -0.5
$ws.Cells.Item(216, 28).Value = -1
$ws.Cells.Item(216, 29).Value = 0.95

# Row 217
$ws.Cells.Item(217, 2).Value = 6571262
$ws.Cells.Item(217, 6).Value = 'CRB'
$ws.Cells.Item(217, 7).Value = 'Ceara'
$ws.Cells.Item(217, 8).Value = 2
$ws.Cells.Item(217, 9).Value = 0
$ws.Cells.Item(217, 10).Value = 'H'
$ws.Cells.Item(217, 11).Value = 2
$ws.Cells.Item(217, 12).Value = 3.1
$ws.Cells.Item(217, 13).Value = 3.7
$ws.Cells.Item(217, 14).Value = 1.95
$ws.Cells.Item(217, 15).Value = 3.25
$ws.Cells.Item(217, 16).Value = 4.333
$ws.Cells.Item(217, 17).Value = -0.5
$ws.Cells.Item(217, 18).Value = 1.925
$ws.Cells.Item(217, 19).Value = 1.875
$ws.Cells.Item(217, 20).Value = 2
$ws.Cells.Item(217, 21).Value = 1.775
$ws.Cells.Item(217, 22).Value = 2.025
$ws.Cells.Item(217, 23).Value = 0.95
$ws.Cells.Item(217, 24).Value = -1
$ws.Cells.Item(217, 25).Value = -1
$ws.Cells.Item(217, 26).Value = 0.925
$ws.Cells.Item(217, 27).Value = -1
$ws.Cells.Item(217, 28).Value = 0
$ws.Cells.Item(217, 29).Value = -0

# Row 227
$ws.Cells.Item(227, 2).Value = 6576739
$ws.Cells.Item(227, 6).Value = 'Gremio Novorizontino'
$ws.Cells.Item(227, 7).Value = 'Tombense MG'
$ws.Cells.Item(227, 8).Value = 0
$ws.Cells.Item(227, 9).Value = 0
$ws.Cells.Item(227, 10).Value = 'D'
$ws.Cells.Item(227, 11).Value = 1.5
$ws.Cells.Item(227, 12).Value = 3.75
$ws.Cells.Item(227, 13).Value = 6.5
$ws.Cells.Item(227, 14).Value = 1.5
$ws.Cells.Item(227, 15).Value = 4
$ws.Cells.Item(227, 16).Value = 7.5
$ws.Cells.Item(227, 17).Value = -1
$ws.Cells.Item(227, 18).Value = 1.825
$ws.Cells.Item(227, 19).Value = 1.975
$ws.Cells.Item(227, 20).Value = 2.5
$ws.Cells.Item(227, 21).Value = 1.925
$ws.Cells.Item(227, 22).Value = 1.875
$ws.Cells.Item(227, 23).Value = -1
$ws.Cells.Item(227, 24).Value = 3
$ws.Cells.Item(227, 25).Value = -1
$ws.Cells.Item(227, 26).Value = -1
$ws.Cells.Item(227, 27).Value = 0.9750000000000001
$ws.Cells.Item(227, 28).Value = -1
$ws.Cells.Item(227, 29).Value = 0.875

# Row 228
$ws.Cells.Item(228, 2).Value = 6576997
$ws.Cells.Item(228, 6).Value = 'Ituano'
$ws.Cells.Item(228, 7).Value = 'CRB'
$ws.Cells.Item(228, 8).Value = 0
$ws.Cells.Item(228, 9).Value = 0
$ws.Cells.Item(228, 10).Value = 'D'
$ws.Cells.Item(228, 11).Value = 2.375
$ws.Cells.Item(228, 12).Value = 3
$ws.Cells.Item(228, 13).Value = 3.25
$ws.Cells.Item(228, 14).Value = 2.15
$ws.Cells.Item(228, 15).Value = 3
$ws.Cells.Item(228, 16).Value = 3.8
$ws.Cells.Item(228, 17).Value = -0.25
$ws.Cells.Item(228, 18).Value = 1.825
$ws.Cells.Item(228, 19).Value = 1.975
$ws.Cells.Item(228, 20).Value = 2
$ws.Cells.Item(228, 21).Value = 1.975
$ws.Cells.Item(228, 22).Value = 1.825
$ws.Cells.Item(228, 23).Value = -1
$ws.Cells.Item(228, 24).Value = 2
$ws.Cells.Item(228, 25).Value = -1
$ws.Cells.Item(228, 26).Value = -0.5
$ws.Cells.Item(228, 27).Value = 0.4875
$ws.Cells.Item(228, 28).Value = -1
$ws.Cells.Item(228, 29).Value = 0.825

# Row 241
$ws.Cells.Item(241, 2).Value = 6590928
$ws.Cells.Item(241, 6).Value = 'Sport Recife'
$ws.Cells.Item(241, 7).Value = 'Chapecoense'
$ws.Cells.Item(241, 8).Value = 2
$ws.Cells.Item(241, 9).Value = 1
$ws.Cells.Item(241, 10).Value = 'H'
$ws.Cells.Item(241, 11).Value = 1.444
$ws.Cells.Item(241, 12).Value = 3.9
$ws.Cells.Item(241, 13).Value = 7
$ws.Cells.Item(241, 14).Value = 1.45
$ws.Cells.Item(241, 15).Value = 4.2
$ws.Cells.Item(241, 16).Value = 7.5
$ws.Cells.Item(241, 17).Value = -1
$ws.Cells.Item(241, 18).Value = 1.75
$ws.Cells.Item(241, 19).Value = 2.05
$ws.Cells.Item(241, 20).Value = 2.25
$ws.Cells.Item(241, 21).Value = 1.825
$ws.Cells.Item(241, 22).Value = 1.975
$ws.Cells.Item(241, 23).Value = 0.45
$ws.Cells.Item(241, 24).Value = -1
$ws.Cells.Item(241, 25).Value = -1
$ws.Cells.Item(241, 26).Value = 0
$ws.Cells.Item(241, 27).Value = -0
$ws.Cells.Item(241, 28).Value = 0.825
$ws.Cells.Item(241, 29).Value = -1

# Row 242
$ws.Cells.Item(242, 2).Value = 6586605
$ws.Cells.Item(242, 6).Value = 'Avai'
$ws.Cells.Item(242, 7).Value = 'Ceara'
$ws.Cells.Item(242, 8).Value = 1
$ws.Cells.Item(242, 9).Value = 0
$ws.Cells.Item(242, 10).Value = 'H'
$ws.Cells.Item(242, 11).Value = 2.4
$ws.Cells.Item(242, 12).Value = 3
$ws.Cells.Item(242, 13).Value = 2.9
$ws.Cells.Item(242, 14).Value = 2.3
$ws.Cells.Item(242, 15).Value = 3.2
$ws.Cells.Item(242, 16).Value = 3.4
$ws.Cells.Item(242, 17).Value = -0.25
$ws.Cells.Item(242, 18).Value = 1.975
$ws.Cells.Item(242, 19).Value = 1.825
$ws.Cells.Item(242, 20).Value = 2
$ws.Cells.Item(242, 21).Value = 1.775
$ws.Cells.Item(242, 22).Value = 2.025
$ws.Cells.Item(242, 23).Value = 1.3
$ws.Cells.Item(242, 24).Value = -1
$ws.Cells.Item(242, 25).Value = -1
$ws.Cells.Item(242, 26).Value = 0.9750000000000001
$ws.Cells.Item(242, 27).Value = -1
$ws.Cells.Item(242, 28).Value = -1
$ws.Cells.Item(242, 29).Value = 1.025

# Row 284
$ws.Cells.Item(284, 2).Value = 6693367
$ws.Cells.Item(284, 6).Value = 'Sport Recife'
$ws.Cells.Item(284, 7).Value = 'Sampaio Correa'
$ws.Cells.Item(284, 8).Value = 4
$ws.Cells.Item(284, 9).Value = 1
$ws.Cells.Item(284, 10).Value = 'H'
$ws.Cells.Item(284, 11).Value = 1.444
$ws.Cells.Item(284, 12).Value = 4
$ws.Cells.Item(284, 13).Value = 7
$ws.Cells.Item(284, 14).Value = 1.533
$ws.Cells.Item(284, 15).Value = 4
$ws.Cells.Item(284, 16).Value = 6
$ws.Cells.Item(284, 17).Value = -1
$ws.Cells.Item(284, 18).Value = 1.975
$ws.Cells.Item(284, 19).Value = 1.825
$ws.Cells.Item(284, 20).Value = 2.25
$ws.Cells.Item(284, 21).Value = 1.775
$ws.Cells.Item(284, 22).Value = 2.025
$ws.Cells.Item(284, 23).Value = 0.5329999999999999
$ws.Cells.Item(284, 24).Value = -1
$ws.Cells.Item(284, 25).Value = -1
$ws.Cells.Item(284, 26).Value = 0.9750000000000001
$ws.Cells.Item(284, 27).Value = -1
$ws.Cells.Item(284, 28).Value = 0.7749999999999999
$ws.Cells.Item(284, 29).Value = -1

# Row 285
$ws.Cells.Item(285, 2).Value = 6693031
$ws.Cells.Item(285, 6).Value = 'Chapecoense'
$ws.Cells.Item(285, 7).Value = 'Vitoria'
$ws.Cells.Item(285, 8).Value = 3
$ws.Cells.Item(285, 9).Value = 1
$ws.Cells.Item(285, 10).Value = 'H'
$ws.Cells.Item(285, 11).Value = 1.8
$ws.Cells.Item(285, 12).Value = 3.6
$ws.Cells.Item(285, 13).Value = 4.333
$ws.Cells.Item(285, 14).Value = 1.615
$ws.Cells.Item(285, 15).Value = 4
$ws.Cells.Item(285, 16).Value = 5.25
$ws.Cells.Item(285, 17).Value = -1
$ws.Cells.Item(285, 18).Value = 2.025
$ws.Cells.Item(285, 19).Value = 1.775
$ws.Cells.Item(285, 20).Value = 2.25
$ws.Cells.Item(285, 21).Value = 1.775
$ws.Cells.Item(285, 22).Value = 2.025
$ws.Cells.Item(285, 23).Value = 0.615
$ws.Cells.Item(285, 24).Value = -1
$ws.Cells.Item(285, 25).Value = -1
$ws.Cells.Item(285, 26).Value = 1.025
$ws.Cells.Item(285, 27).Value = -1
$ws.Cells.Item(285, 28).Value = 0.7749999999999999
$ws.Cells.Item(285, 29).Value = -1

# Row 286
$ws.Cells.Item(286, 2).Value = 6693030
$ws.Cells.Item(286, 6).Value = 'Ceara'
$ws.Cells.Item(286, 7).Value = 'EC Juventude'
$ws.Cells.Item(286, 8).Value = 1
$ws.Cells.Item(286, 9).Value = 3
$ws.Cells.Item(286, 10).Value = 'A'
$ws.Cells.Item(286, 11).Value = 3.25
$ws.Cells.Item(286, 12).Value = 3.4
$ws.Cells.Item(286, 13).Value = 2.1
$ws.Cells.Item(286, 14).Value = 3.3
$ws.Cells.Item(286, 15).Value = 3.2
$ws.Cells.Item(286, 16).Value = 2.3
$ws.Cells.Item(286, 17).Value = 0.25
$ws.Cells.Item(286, 18).Value = 1.875
$ws.Cells.Item(286, 19).Value = 1.975
$ws.Cells.Item(286, 20).Value = 2
$ws.Cells.Item(286, 21).Value = 1.825
$ws.Cells.Item(286, 22).Value = 2.025
$ws.Cells.Item(286, 23).Value = -1
$ws.Cells.Item(286, 24).Value = -1
$ws.Cells.Item(286, 25).Value = 1.3
$ws.Cells.Item(286, 26).Value = -1
$ws.Cells.Item(286, 27).Value = 0.9750000000000001
$ws.Cells.Item(286, 28).Value = 0.825
$ws.Cells.Item(286, 29).Value = -1

# Row 287
$ws.Cells.Item(287, 2).Value = 6693028
$ws.Cells.Item(287, 6).Value = 'Ponte Preta'
$ws.Cells.Item(287, 7).Value = 'CRB'
$ws.Cells.Item(287, 8).Value = 3
$ws.Cells.Item(287, 9).Value = 0
$ws.Cells.Item(287, 10).Value = 'H'
$ws.Cells.Item(287, 11).Value = 1.727
$ws.Cells.Item(287, 12).Value = 3.5
$ws.Cells.Item(287, 13).Value = 4
$ws.Cells.Item(287, 14).Value = 1.7
$ws.Cells.Item(287, 15).Value = 3.6
$ws.Cells.Item(287, 16).Value = 5
$ws.Cells.Item(287, 17).Value = -0.75
$ws.Cells.Item(287, 18).Value = 1.975
$ws.Cells.Item(287, 19).Value = 1.875
$ws.Cells.Item(287, 20).Value = 2
$ws.Cells.Item(287, 21).Value = 1.775
$ws.Cells.Item(287, 22).Value = 2.1
$ws.Cells.Item(287, 23).Value = 0.7
$ws.Cells.Item(287, 24).Value = -1
$ws.Cells.Item(287, 25).Value = -1
$ws.Cells.Item(287, 26).Value = 0.9750000000000001
$ws.Cells.Item(287, 27).Value = -1
$ws.Cells.Item(287, 28).Value = 0.7749999999999999
$ws.Cells.Item(287, 29).Value = -1

# Row 289
$ws.Cells.Item(289, 2).Value = 6689427
$ws.Cells.Item(289, 6).Value = 'Avai'
$ws.Cells.Item(289, 7).Value = 'Ituano'
$ws.Cells.Item(289, 8).Value = 0
$ws.Cells.Item(289, 9).Value = 0
$ws.Cells.Item(289, 10).Value = 'D'
$ws.Cells.Item(289, 11).Value = 1.95
$ws.Cells.Item(289, 12).Value = 3.1
$ws.Cells.Item(289, 13).Value = 4.2
$ws.Cells.Item(289, 14).Value = 2.4
$ws.Cells.Item(289, 15).Value = 3.2
$ws.Cells.Item(289, 16).Value = 3.2
$ws.Cells.Item(289, 17).Value = -0.25
$ws.Cells.Item(289, 18).Value = 2.05
$ws.Cells.Item(289, 19).Value = 1.75
$ws.Cells.Item(289, 20).Value = 2.25
$ws.Cells.Item(289, 21).Value = 1.825
$ws.Cells.Item(289, 22).Value = 1.975
$ws.Cells.Item(289, 23).Value = -1
$ws.Cells.Item(289, 24).Value = 2.2
$ws.Cells.Item(289, 25).Value = -1
$ws.Cells.Item(289, 26).Value = -0.5
$ws.Cells.Item(289, 27).Value = 0.375
$ws.Cells.Item(289, 28).Value = -1
$ws.Cells.Item(289, 29).Value = 0.9750000000000001

# Row 290
$ws.Cells.Item(290, 2).Value = 6693029
$ws.Cells.Item(290, 6).Value = 'Atletico GO'
$ws.Cells.Item(290, 7).Value = 'Guarani'
$ws.Cells.Item(290, 8).Value = 3
$ws.Cells.Item(290, 9).Value = 0
$ws.Cells.Item(290, 10).Value = 'H'
$ws.Cells.Item(290, 11).Value = 1.45
$ws.Cells.Item(290, 12).Value = 4
$ws.Cells.Item(290, 13).Value = 8
$ws.Cells.Item(290, 14).Value = 1.333
$ws.Cells.Item(290, 15).Value = 4.8
$ws.Cells.Item(290, 16).Value = 10
$ws.Cells.Item(290, 17).Value = -1.25
$ws.Cells.Item(290, 18).Value = 1.8
$ws.Cells.Item(290, 19).Value = 2
$ws.Cells.Item(290, 20).Value = 2.5
$ws.Cells.Item(290, 21).Value = 1.95
$ws.Cells.Item(290, 22).Value = 1.85
$ws.Cells.Item(290, 23).Value = 0.333
$ws.Cells.Item(290, 24).Value = -1
$ws.Cells.Item(290, 25).Value = -1
$ws.Cells.Item(290, 26).Value = 0.8
$ws.Cells.Item(290, 27).Value = -1
$ws.Cells.Item(290, 28).Value = 0.95
$ws.Cells.Item(290, 29).Value = -1

# Row 291
$ws.Cells.Item(291, 2).Value = 6689350
$ws.Cells.Item(291, 6).Value = 'Tombense MG'
$ws.Cells.Item(291, 7).Value = 'Mirassol'
$ws.Cells.Item(291, 8).Value = 0
$ws.Cells.Item(291, 9).Value = 1
$ws.Cells.Item(291, 10).Value = 'A'
$ws.Cells.Item(291, 11).Value = 3.2
$ws.Cells.Item(291, 12).Value = 3
$ws.Cells.Item(291, 13).Value = 2.25
$ws.Cells.Item(291, 14).Value = 3
$ws.Cells.Item(291, 15).Value = 3.25
$ws.Cells.Item(291, 16).Value = 2.3
$ws.Cells.Item(291, 17).Value = 0.25
$ws.Cells.Item(291, 18).Value = 1.775
$ws.Cells.Item(291, 19).Value = 2.025
$ws.Cells.Item(291, 20).Value = 2.5
$ws.Cells.Item(291, 21).Value = 2
$ws.Cells.Item(291, 22).Value = 1.8
$ws.Cells.Item(291, 23).Value = -1
$ws.Cells.Item(291, 24).Value = -1
$ws.Cells.Item(291, 25).Value = 1.3
$ws.Cells.Item(291, 26).Value = -1
$ws.Cells.Item(291, 27).Value = 1.025
$ws.Cells.Item(291, 28).Value = -1
$ws.Cells.Item(291, 29).Value = 0.8

# Row 292
$ws.Cells.Item(292, 2).Value = 6689425
$ws.Cells.Item(292, 6).Value = 'Gremio Novorizontino'
$ws.Cells.Item(292, 7).Value = 'Criciuma'
$ws.Cells.Item(292, 8).Value = 2
$ws.Cells.Item(292, 9).Value = 0
$ws.Cells.Item(292, 10).Value = 'H'
$ws.Cells.Item(292, 11).Value = 1.571
$ws.Cells.Item(292, 12).Value = 3.8
$ws.Cells.Item(292, 13).Value = 5.75
$ws.Cells.Item(292, 14).Value = 1.45
$ws.Cells.Item(292, 15).Value = 4.5
$ws.Cells.Item(292, 16).Value = 6.5
$ws.Cells.Item(292, 17).Value = -1
$ws.Cells.Item(292, 18).Value = 1.75
$ws.Cells.Item(292, 19).Value = 2.05
$ws.Cells.Item(292, 20).Value = 2.5
$ws.Cells.Item(292, 21).Value = 1.975
$ws.Cells.Item(292, 22).Value = 1.825
$ws.Cells.Item(292, 23).Value = 0.45
$ws.Cells.Item(292, 24).Value = -1
$ws.Cells.Item(292, 25).Value = -1
$ws.Cells.Item(292, 26).Value = 0.75
$ws.Cells.Item(292, 27).Value = -1
$ws.Cells.Item(292, 28).Value = -1
$ws.Cells.Item(292, 29).Value = 0.825

# Row 293
$ws.Cells.Item(293, 2).Value = 6689429
$ws.Cells.Item(293, 6).Value = 'ABC'
$ws.Cells.Item(293, 7).Value = 'Vila Nova'
$ws.Cells.Item(293, 8).Value = 3
$ws.Cells.Item(293, 9).Value = 2
$ws.Cells.Item(293, 10).Value = 'H'
$ws.Cells.Item(293, 11).Value = 8
$ws.Cells.Item(293, 12).Value = 4.75
$ws.Cells.Item(293, 13).Value = 1.363
$ws.Cells.Item(293, 14).Value = 6.5
$ws.Cells.Item(293, 15).Value = 4.2
$ws.Cells.Item(293, 16).Value = 1.45
$ws.Cells.Item(293, 17).Value = 1
$ws.Cells.Item(293, 18).Value = 2
$ws.Cells.Item(293, 19).Value = 1.8
$ws.Cells.Item(293, 20).Value = 2.25
$ws.Cells.Item(293, 21).Value = 1.9
$ws.Cells.Item(293, 22).Value = 1.9
$ws.Cells.Item(293, 23).Value = 5.5
$ws.Cells.Item(293, 24).Value = -1
$ws.Cells.Item(293, 25).Value = -1
$ws.Cells.Item(293, 26).Value = 1
$ws.Cells.Item(293, 27).Value = -1
$ws.Cells.Item(293, 28).Value = 0.8999999999999999
$ws.Cells.Item(293, 29).Value = -1
